$d = $word.ActiveDocument

# --- Part 1: remove the stray "_GoBack" bookmark that wraps the DOMAIN MODEL
# image (it is re-created later at the real last-edit location, see Part 2).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Part 2: the "3) The system sends notification about acceptance or
# rejection." paragraph loses its text, and the following empty (bold)
# paragraph is removed entirely, merging into a single empty paragraph
# that keeps the first paragraph's properties. A fresh "_GoBack" bookmark
# is left behind in that now-empty paragraph, marking the last edit.
$rng = $d.Content
$found = $rng.Find.Execute("3) The system sends notification about acceptance or rejection.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Delete()

    $para = $rng.Paragraphs(1)
    $nextPara = $para.Next()
    $nextPara.Range.Delete()

    $d.Bookmarks.Add("_GoBack", $rng)
}
